$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = "excel11"
$ws.Range("B3").Value = "excel12"
$ws.Range("B4").Value = "excel33"
$ws.Range("B5").Value = "excel44"
$ws.Range("D9").Select()
